$wb = $excel.ActiveWorkbook

# Update the selection on the existing RegisterData sheet (E9 -> D9)
$ws1 = $wb.Worksheets.Item("RegisterData")
[void]$ws1.Range("D9").Select()

# Add the new "DeleteRecord" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "DeleteRecord"

# Populate the new sheet with data-driven values for the DeleteRecord test
$ws2.Range("A1").Value = "UserName  "
$ws2.Range("A2").Value = "novak"
$ws2.Columns.Item(1).AutoFit()

[void]$ws2.Range("A4").Select()
